$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 122.5
$ws.Range("I9").Value = 87.09999999999999
$ws.Range("J9").Value = 299.5
$ws.Range("K9").Value = 87.09999999999999
$ws.Range("L9").Value = 299.5
$ws.Range("M9").Value = 81.90000000000001
$ws.Range("N9").Value = -637.5
$ws.Range("H47").Value = 6783.5
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H53").Value = 417.2857
$ws.Range("I53").Value = 320.16666
$ws.Range("K53").Value = 320.16666
$ws.Range("M53").Value = 316.83334
$ws.Range("H92").Value = 181.66667
$ws.Range("I92").Value = 181.66667
$ws.Range("K92").Value = 181.66667
$ws.Range("M92").Value = 1066.33333
$ws.Range("H99").Value = 5364.8335
$ws.Range("J99").Value = 5836.4
$ws.Range("L99").Value = 17509.2
$ws.Range("N99").Value = -20505.2
$ws.Range("H104").Value = 446.5
$ws.Range("I104").Value = 455.8
$ws.Range("J104").Value = 400
$ws.Range("K104").Value = 1367.4
$ws.Range("L104").Value = 1200
$ws.Range("M104").Value = 379.5999999999999
$ws.Range("N104").Value = -4694
$ws.Range("H116").Value = 3133.625
$ws.Range("I116").Value = 3093.3333
$ws.Range("K116").Value = 3093.3333
$ws.Range("M116").Value = 348.6667000000002
$ws.Range("H135").Value = 916.6667
$ws.Range("I135").Value = 916.6667
$ws.Range("K135").Value = 8250.0003
$ws.Range("M135").Value = -5715.0003
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 891.7273
$ws.Range("I2").Value = 891.7273
$ws.Range("K2").Value = 891.7273
$ws.Range("M2").Value = -778.7273
$ws.Range("H32").Value = 3281.3333
$ws.Range("I32").Value = 3281.3333
$ws.Range("K32").Value = 3281.3333
$ws.Range("M32").Value = -2994.3333
$ws.Range("H45").Value = 2407.2222
$ws.Range("I45").Value = 1410
$ws.Range("J45").Value = 3653.75
$ws.Range("K45").Value = 1410
$ws.Range("L45").Value = 3653.75
$ws.Range("M45").Value = -1033
$ws.Range("N45").Value = -4407.75
$ws.Range("H61").Value = 4136.9
$ws.Range("I61").Value = 3278.7058
$ws.Range("K61").Value = 3278.7058
$ws.Range("M61").Value = -3066.7058
$ws.Range("H69").Value = 249999
$ws.Range("J69").Value = 249999
$ws.Range("L69").Value = 249999
$ws.Range("N69").Value = -251497
$ws.Range("H72").Value = 249999
$ws.Range("J72").Value = 249999
$ws.Range("L72").Value = 749997
$ws.Range("N72").Value = -757485
$ws.Range("H116").Value = 891.7273
$ws.Range("I116").Value = 891.7273
$ws.Range("K116").Value = 891.7273
$ws.Range("M116").Value = 1402.2727
$ws.Range("H122").Value = 716.6667
$ws.Range("I122").Value = 716.6667
$ws.Range("K122").Value = 2150.0001
$ws.Range("M122").Value = 299.9998999999998
$ws.Range("H132").Value = 7193.8887
$ws.Range("I132").Value = 7193.8887
$ws.Range("K132").Value = 21581.6661
$ws.Range("M132").Value = -19051.6661
$ws.Range("H136").Value = 4136.9
$ws.Range("I136").Value = 3278.7058
$ws.Range("K136").Value = 9836.117400000001
$ws.Range("M136").Value = -7286.117400000001
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 891.7273
$ws.Range("I3").Value = 891.7273
$ws.Range("K3").Value = 891.7273
$ws.Range("M3").Value = -777.7273
$ws.Range("H26").Value = 16980.125
$ws.Range("I26").Value = 16980.125
$ws.Range("K26").Value = 16980.125
$ws.Range("M26").Value = -16688.125
$ws.Range("H86").Value = 3937.4285
$ws.Range("I86").Value = 1798.9231
$ws.Range("K86").Value = 1798.9231
$ws.Range("M86").Value = -675.9231
$ws.Range("H89").Value = 3937.4285
$ws.Range("I89").Value = 1798.9231
$ws.Range("K89").Value = 8994.6155
$ws.Range("M89").Value = -3378.6155
$ws.Range("H96").Value = 14450.875
$ws.Range("I96").Value = 14450.875
$ws.Range("K96").Value = 14450.875
$ws.Range("M96").Value = -11704.875
$ws.Range("H105").Value = 1567.2858
$ws.Range("I105").Value = 1387.6364
$ws.Range("K105").Value = 1387.6364
$ws.Range("M105").Value = 359.3635999999999
$ws.Range("H134").Value = 1464.6
$ws.Range("I134").Value = 837
$ws.Range("J134").Value = 3975
$ws.Range("K134").Value = 2511
$ws.Range("L134").Value = 11925
$ws.Range("M134").Value = 24
$ws.Range("N134").Value = -16995
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 55000
$ws.Range("J97").Value = 55000
$ws.Range("L97").Value = 55000
$ws.Range("N97").Value = -56982
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 295
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 295
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 885
$ws.Range("N22").Value = -1223
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 295
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 295
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 885
$ws.Range("N27").Value = -1089
$ws.Range("M27").ClearContents()
$ws.Range("H42").Value = 2000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 2000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 6000
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -7068
$ws.Range("H68").Value = 1058.6666
$ws.Range("I68").Value = 926
$ws.Range("J68").Value = 1125
$ws.Range("K68").Value = 2778
$ws.Range("L68").Value = 3375
$ws.Range("M68").Value = -1967
$ws.Range("N68").Value = -4997
$ws.Range("H71").Value = 1058.6666
$ws.Range("I71").Value = 926
$ws.Range("J71").Value = 1125
$ws.Range("K71").Value = 8334
$ws.Range("L71").Value = 10125
$ws.Range("M71").Value = -4278
$ws.Range("N71").Value = -18237
$ws.Range("H115").Value = 1995
$ws.Range("I115").Value = 1995
$ws.Range("K115").Value = 5985
$ws.Range("M115").Value = -4810
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2033.7778
$ws.Range("I80").Value = 1783
$ws.Range("K80").Value = 1783
$ws.Range("M80").Value = -785
$ws.Range("H83").Value = 2033.7778
$ws.Range("I83").Value = 1783
$ws.Range("K83").Value = 8915
$ws.Range("M83").Value = -3923
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2901.818
$ws.Range("I122").Value = 2634
$ws.Range("J122").Value = 3002.25
$ws.Range("K122").Value = 7902
$ws.Range("L122").Value = 9006.75
$ws.Range("M122").Value = -5452
$ws.Range("N122").Value = -13906.75
$ws.Range("H136").Value = 3049.1667
$ws.Range("I136").Value = 2849.5
$ws.Range("J136").Value = 3149
$ws.Range("K136").Value = 8548.5
$ws.Range("L136").Value = 9447
$ws.Range("M136").Value = -5998.5
$ws.Range("N136").Value = -14547
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 34774.75
$ws.Range("J41").Value = 34494.668
$ws.Range("L41").Value = 34494.668
$ws.Range("N41").Value = -35274.668
$ws.Range("H81").Value = 3000
$ws.Range("I81").Value = 3000
$ws.Range("K81").Value = 6000
$ws.Range("M81").Value = -4939
$ws.Range("H84").Value = 3000
$ws.Range("I84").Value = 3000
$ws.Range("K84").Value = 30000
$ws.Range("M84").Value = -24696
$ws.Range("H93").Value = 33389
$ws.Range("J93").Value = 33389
$ws.Range("L93").Value = 33389
$ws.Range("N93").Value = -38381
$ws.Range("H122").Value = 2494
$ws.Range("I122").Value = 1376.8462
$ws.Range("K122").Value = 4130.5386
$ws.Range("M122").Value = -1680.5386
$ws.Range("H132").Value = 2315
$ws.Range("I132").Value = 2086.8333
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 6260.499899999999
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -3730.499899999999
$ws.Range("N132").Value = -14058.5
